$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.441.85'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.663.29'
$ws.Range('E3').Value = '  +1.21%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9983'
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '235.98'
$ws.Range('E5').Value = '  +1.12%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4624'
$ws.Range('E7').Value = '  -3.15%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2571'
$ws.Range('E8').Value = '  -1.32%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06140'
$ws.Range('E9').Value = '  +0.51%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.657.89'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06954'
$ws.Range('E11').Value = '  -1.30%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.58'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.319'
$ws.Range('E13').Value = '  -0.79%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '75.02'
$ws.Range('E14').Value = '  +1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.5716'
$ws.Range('E15').Value = '  -4.05%  '
$ws.Range('E16').Value = '  +0.13%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.000'
$ws.Range('E17').Value = '  +0.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.366.80'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000006686'
$ws.Range('E19').Value = '  +1.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.33'
$ws.Range('E20').Value = '  +0.51%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.873.38'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.397'
$ws.Range('E22').Value = '  +1.71%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.603'
$ws.Range('E23').Value = '  +0.39%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.208'
$ws.Range('E24').Value = '  -0.73%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '134.24'
$ws.Range('E25').Value = '  -0.34%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '14.88'
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.365'
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.709'
$ws.Range('E28').Value = '  +4.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '103.89'
$ws.Range('E29').Value = '  +0.14%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.930'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07673'
$ws.Range('E31').Value = '  -0.31%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.582'
$ws.Range('E32').Value = '  +0.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04339'
$ws.Range('E33').Value = '  +1.09%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.629'
$ws.Range('E34').Value = '  +2.22%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6003'
$ws.Range('E35').Value = '  +1.08%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.9360'
$ws.Range('E36').Value = '  +0.92%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9116'
$ws.Range('E37').Value = '  +2.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '107.28'
$ws.Range('E38').Value = '  +8.41%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9988'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.230'
$ws.Range('E40').Value = '  -13.72%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.819'
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('E42').Value = '  -4.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.3692'
$ws.Range('E43').Value = '  -0.39%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.986'
$ws.Range('E44').Value = '  +6.65%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1103'
$ws.Range('E45').Value = '  +0.12%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.05255'
$ws.Range('E46').Value = '  +0.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.092'
$ws.Range('E47').Value = '  -0.09%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.33'
$ws.Range('E48').Value = '  +4.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.519'
$ws.Range('E49').Value = '  +5.64%  '
$ws.Range('B50').Value = 'TrueUSD'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1ZZI6g5k5royD+trueusd-tusd'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.002'
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9975'
$ws.Range('E51').Value = '  -0.11%  '
